$d = $word.ActiveDocument

# Locate the (currently blank) "Liste points" bullet paragraph that sits
# right under the "Liste des points de matières" heading and right above
# the "Application" Titre1 heading.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Liste points") {
        $trimmed = $p.Range.Text.Trim()
        if ($trimmed -eq "") {
            $target = $p
            break
        }
    }
}

if ($target -eq $null) {
    Write-Host "ERROR: placeholder bullet paragraph not found"
} else {
    # Turn the blank bullet into the first new subject-matter point.
    $r = $target.Range
    $r.Find.Execute(" ", $false, $false, $false, $false, $false, $true, 1, $false, "Pattern : modèle vue contrôleur.", 2)

    # Re-fetch the (now filled-in) paragraph and add a second bullet point
    # right after it, matching the same "Liste points" style.
    $first = $target
    $first.Range.InsertParagraphAfter()

    $second = $first.Next()
    $second.Range.Text = "Utilisation de Git"
}
